# Insert two new weekly price records for "Poroto verde" (Brío / Sin especificar)
# at the top of the historical block (rows 88-89), shifting the existing
# rows 88..179 down to 90..181.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88:A89").EntireRow.Insert()

# Row 88: new "Brío" record
$ws.Cells.Item(88,1).Value2  = 10
$ws.Cells.Item(88,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(88,3).Value2  = "La Araucanía"
$ws.Cells.Item(88,4).Value2  = 44966
$ws.Cells.Item(88,5).Value2  = 9
$ws.Cells.Item(88,6).Value2  = 100112031
$ws.Cells.Item(88,7).Value2  = "Poroto verde"
$ws.Cells.Item(88,8).Value2  = "Brío"
$ws.Cells.Item(88,9).Value2  = "Primera"
$ws.Cells.Item(88,10).Value2 = 200
$ws.Cells.Item(88,11).Value2 = 1600
$ws.Cells.Item(88,12).Value2 = 1600
$ws.Cells.Item(88,13).Value2 = 1600
$ws.Cells.Item(88,14).Value2 = "$/kilo"
$ws.Cells.Item(88,15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(88,16).Value2 = 1600
$ws.Cells.Item(88,17).Value2 = 1
$ws.Cells.Item(88,18).Value2 = "Hortaliza"

# Row 89: duplicate new record
$ws.Cells.Item(89,1).Value2  = 10
$ws.Cells.Item(89,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(89,3).Value2  = "La Araucanía"
$ws.Cells.Item(89,4).Value2  = 44966
$ws.Cells.Item(89,5).Value2  = 9
$ws.Cells.Item(89,6).Value2  = 100112031
$ws.Cells.Item(89,7).Value2  = "Poroto verde"
$ws.Cells.Item(89,8).Value2  = "Brío"
$ws.Cells.Item(89,9).Value2  = "Primera"
$ws.Cells.Item(89,10).Value2 = 200
$ws.Cells.Item(89,11).Value2 = 1600
$ws.Cells.Item(89,12).Value2 = 1600
$ws.Cells.Item(89,13).Value2 = 1600
$ws.Cells.Item(89,14).Value2 = "$/kilo"
$ws.Cells.Item(89,15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(89,16).Value2 = 1600
$ws.Cells.Item(89,17).Value2 = 1
$ws.Cells.Item(89,18).Value2 = "Hortaliza"
